$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-9 ---

# Row 2: Literature Society IITJ Website
$ws.Cells.Item(2,3).Value = 82.95
$ws.Cells.Item(2,5).Value = 74.66

# Row 3: name changes to LLMGuard, scores updated
$ws.Cells.Item(3,2).Value = "LLMGuard"
$ws.Cells.Item(3,3).Value = 76.23
$ws.Cells.Item(3,5).Value = 76.23

# Row 4: CloudPhysician's Vital Extraction Challenge
$ws.Cells.Item(4,3).Value = 69.75
$ws.Cells.Item(4,5).Value = 59.29

# Row 5: SMART SENSING MIDDLEWARE
$ws.Cells.Item(5,3).Value = 108.75
$ws.Cells.Item(5,5).Value = 100

# Row 6: RAPID
$ws.Cells.Item(6,3).Value = 108.75
$ws.Cells.Item(6,5).Value = 100

# Row 7: SHAMIYANA APP
$ws.Cells.Item(7,3).Value = 79.56
$ws.Cells.Item(7,5).Value = 71.59999999999999

# Row 8: Website for the Literature Society of the college
$ws.Cells.Item(8,3).Value = 82.95
$ws.Cells.Item(8,5).Value = 74.66
$ws.Cells.Item(8,6).Value = 3

# Row 9: LLMGuard (name unchanged)
$ws.Cells.Item(9,3).Value = 79.56
$ws.Cells.Item(9,5).Value = 79.56
$ws.Cells.Item(9,6).Value = 3

# --- Add new rows 10-13 ---

# Row 10
$ws.Cells.Item(10,1).Value = 3
$ws.Cells.Item(10,2).Value = "Multi Model Data Analysis for Annotation of Human Activities"
$ws.Cells.Item(10,3).Value = 76.23
$ws.Cells.Item(10,4).Value = 1
$ws.Cells.Item(10,5).Value = 76.23
$ws.Cells.Item(10,6).Value = 3

# Row 11
$ws.Cells.Item(11,1).Value = 4
$ws.Cells.Item(11,2).Value = "Video Conferencing Project"
$ws.Cells.Item(11,3).Value = 86.40000000000001
$ws.Cells.Item(11,4).Value = 0.85
$ws.Cells.Item(11,5).Value = 73.44
$ws.Cells.Item(11,6).Value = 3

# Row 12
$ws.Cells.Item(12,1).Value = 4
$ws.Cells.Item(12,2).Value = "Alcheringa Pass Portal"
$ws.Cells.Item(12,3).Value = 86.40000000000001
$ws.Cells.Item(12,4).Value = 1
$ws.Cells.Item(12,5).Value = 86.40000000000001
$ws.Cells.Item(12,6).Value = 3

# Row 13
$ws.Cells.Item(13,1).Value = 4
$ws.Cells.Item(13,2).Value = "TEDxIITGuwahati Website"
$ws.Cells.Item(13,3).Value = 79.56
$ws.Cells.Item(13,4).Value = 0.85
$ws.Cells.Item(13,5).Value = 67.63
$ws.Cells.Item(13,6).Value = 3
